$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" property to the new generation timestamp.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2023-11-27T16:21:02+00:00"

# ---------------------------------------------------------------------------
# 2. Concepts sheet: add Somatic counterparts next to each Germline concept.
#    Row layout before edit (row -> Level/Code/Display):
#       2  ALIR
#       3  SNV  (Display was placeholder "SNV", now becomes "Germline SNV")
#       4  GCNV / Germline CNV
#       5  GSV  / Germline SV
#       6..14 remaining concepts
#
#    New layout after edit:
#       2  ALIR
#       3  SNV   / Germline SNV
#       4  SSNV  / Somatic SNV      <- new
#       5  GCNV  / Germline CNV
#       6  SCNV  / Somatic CNV      <- new
#       7  GSV   / Germline SV
#       8  SSV   / Somatic SV       <- new
#       9..17 remaining concepts (shifted down by 3)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Concepts")

# Insert the three new blank rows from the bottom up so earlier row numbers
# used below stay valid while we work.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(4).Insert()

# Copy the formatting (border/fill/font/alignment) of an existing data row
# onto the freshly inserted rows so they render identically to the rest of
# the table.
$ws.Range("A2:D2").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Range("A8:D8").PasteSpecial(-4122)

# Column A is always the literal text "1" (the concept hierarchy level).
# Copy it (value + format) from an existing row so it stays a text value
# rather than being re-interpreted as a number.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4104)
$ws.Range("A6").PasteSpecial(-4104)
$ws.Range("A8").PasteSpecial(-4104)

# Correct the existing SNV row's Display column.
$ws.Cells.Item(3, 3).Value = "Germline SNV"

# Populate the new Somatic SNV row.
$ws.Cells.Item(4, 2).Value = "SSNV"
$ws.Cells.Item(4, 3).Value = "Somatic SNV"

# Populate the new Somatic CNV row.
$ws.Cells.Item(6, 2).Value = "SCNV"
$ws.Cells.Item(6, 3).Value = "Somatic CNV"

# Populate the new Somatic SV row.
$ws.Cells.Item(8, 2).Value = "SSV"
$ws.Cells.Item(8, 3).Value = "Somatic SV"
